# Finish Lesson 5 Quiz
# Log 0.5 hours for JS101 on the row for 11/21/2021 (row 89), with
# the Milestones note "Finish Lesson 5" and Notes "Finish Lesson 5 Quiz".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B89").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C89").Value = 0.5

# Set E89 (Milestones) before D89 (Notes) so the new shared strings are
# appended to the shared-string table in the same order as the source edit.
$ws.Range("E89").Value = "Finish Lesson 5"
$ws.Range("D89").Value = "Finish Lesson 5 Quiz"

# Move the active selection to D90, matching the cursor position left by
# the author after making the edit.
$ws.Range("D90").Select() | Out-Null
